$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------------
# 1. D37: turn on "Wrap Text" (keeps its green fill/value, just adds wrap)
# ---------------------------------------------------------------------------
$ws.Range("D37").WrapText = $true

# ---------------------------------------------------------------------------
# 2. Add a new row 38 (duplicate of the "show-add-data" request that used to
#    live on row 42, now documented right under row 37's "show-edit-data").
# ---------------------------------------------------------------------------
$ws.Range("A38").Value = $ws.Range("A37").Value2
$ws.Range("B38").Value = $ws.Range("B37").Value2
$ws.Range("C38").Value = $ws.Range("C37").Value2

# D38 re-uses the plain green-fill style (same style D37 had before the wrap
# text was turned on) -- copy it from D36, which still carries that style.
$ws.Range("D36").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D38").Value = "show-add-data"

$ws.Range("E38").Value = "a compléter (l'url est renvoyée depuis le serveur)"
$ws.Range("E38").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 3. Rows 41-43: the obsolete "ajaxgetchildren" / "show-add-data" /
#    "proxy"+"getfeatureinfo" request rows are cleaned out (content removed,
#    formatting on D kept).
# ---------------------------------------------------------------------------
$ws.Range("C41").ClearContents()
$ws.Range("D41").ClearContents()
$ws.Range("C42").ClearContents()
$ws.Range("D42").ClearContents()
$ws.Range("C43").ClearContents()
$ws.Range("D43").ClearContents()

# ---------------------------------------------------------------------------
# 4. Update the view: scroll down and select C41:D41.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C41:D41").Select()
